$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the activity log text for row 18 (F18) - append the new sentence.
$ws.Range("F18").Value = 'Added new "Likes" table wrote backend code to update likes in the database. Added backend code to close a session.. Fixed bugs in deleting message and fixed bugs in Junit tests involving the test to update message content and the test to delete messages. Fixed a JUnit test bug where we were testing the wrong procedure. Wrote SPROC to toggle likes on a message. Complete close session. Wrote new Junit tests for bad delete requests.'

# Update the hours value for row 18 (E18) from 5 to 6.
$ws.Range("E18").Value = 6

# Update the selected cell to E19 (matches the active cell moved after editing).
$ws.Range("E19").Select()
